$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.032.19"
$ws.Range("E2").Value = "  +2.81%  "

$ws.Range("D3").Value = "3.805.11"
$ws.Range("E3").Value = "  +0.83%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").Value = "'708.58"
$ws.Range("E5").Value = "  +12.55%  "

$ws.Range("D6").Value = "'173.42"
$ws.Range("E6").Value = "  +4.85%  "

$ws.Range("D7").Value = "3.802.20"
$ws.Range("E7").Value = "  +0.83%  "

$ws.Range("E8").Value = "  -0.04%  "

$ws.Range("E9").Value = "  +1.47%  "

$ws.Range("E10").Value = "  +3.39%  "

$ws.Range("E11").Value = "  +9.53%  "

$ws.Range("E12").Value = "  +1.59%  "

$ws.Range("E13").Value = "  +9.63%  "

$ws.Range("D14").Value = "'36.38"
$ws.Range("E14").Value = "  +4.40%  "

$ws.Range("D15").Value = "4.444.96"
$ws.Range("E15").Value = "  +0.89%  "

$ws.Range("D16").Value = "3.806.03"
$ws.Range("E16").Value = "  +0.83%  "

$ws.Range("D17").Value = "71.063.60"
$ws.Range("E17").Value = "  +2.90%  "

$ws.Range("D18").Value = "'17.90"
$ws.Range("E18").Value = "  +1.38%  "

$ws.Range("E19").Value = "  +3.25%  "

$ws.Range("E20").Value = "  +0.44%  "

$ws.Range("D21").Value = "'11.19"
$ws.Range("E21").Value = "  +17.65%  "

$ws.Range("D22").Value = "'484.47"
$ws.Range("E22").Value = "  +3.38%  "

$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").Value = "'83.94"
$ws.Range("E24").Value = "  +2.30%  "

$ws.Range("B25").Value = "PEPE"
$ws.Range("C25").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D25").Value = "'0.0000147"
$ws.Range("E25").Value = "  +3.49%  "

$ws.Range("E26").Value = "  +2.81%  "

$ws.Range("E27").Value = "  +4.25%  "

$ws.Range("E28").Value = "  +3.04%  "

$ws.Range("D29").Value = "3.955.57"
$ws.Range("E29").Value = "  +0.91%  "

$ws.Range("B30").Value = "Dai"
$ws.Range("C30").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D30").Value = "'0.999"
$ws.Range("E30").Value = "  -0.08%  "

$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").Value = "'3.11"
$ws.Range("E31").Value = "  +16.91%  "

$ws.Range("D32").Value = "'7.61"
$ws.Range("E32").Value = "  +7.11%  "

$ws.Range("E33").Value = "  +0.37%  "

$ws.Range("D34").Value = "'29.68"

$ws.Range("E35").Value = "  +1.12%  "

$ws.Range("D36").Value = "'9.27"
$ws.Range("E36").Value = "  +4.47%  "

$ws.Range("E37").Value = "  +0.04%  "

$ws.Range("D38").Value = "3.755.25"
$ws.Range("E38").Value = "  +0.84%  "

$ws.Range("E39").Value = "  +2.62%  "

$ws.Range("E40").Value = "  +9.54%  "

$ws.Range("E41").Value = "  +3.58%  "

$ws.Range("D42").Value = "'2.24"
$ws.Range("E42").Value = "  +12.92%  "

$ws.Range("E43").Value = "  +26.13%  "

$ws.Range("D44").Value = "'0.971"

$ws.Range("E45").Value = "  +0.09%  "

$ws.Range("E46").Value = "  -0.02%  "

$ws.Range("D47").Value = "'162.65"
$ws.Range("E47").Value = "  +4.15%  "

$ws.Range("D48").Value = "'49.47"
$ws.Range("E48").Value = "  +5.36%  "

$ws.Range("D49").Value = "'45.06"
$ws.Range("E49").Value = "  +2.87%  "

$ws.Range("E50").Value = "  +2.59%  "

$ws.Range("D51").Value = "'1.37"
$ws.Range("E51").Value = "  -2.23%  "
